# Team19 poster edits:
#  - Subtitle tagline: replace placeholder tagline with the real one
#    (typed as two runs, matching the way PowerPoint records "type to replace
#    selection" + "continue typing" as separate runs).
#  - Client/Team block: "Client(s):" -> "Client:" and drop the "Nasa JPL"
#    affiliation from Kim Whitehall's line.
#  - Big poster number placeholder: "###" -> "CS19".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# --- Subtitle 2 (shape id 13): tagline text ---------------------------------
$subtitle = Get-ShapeById $s 13
$tagline = $subtitle.TextFrame.TextRange
$tagline.Text = "Staring at particles of dust never looked "
$null = $tagline.InsertAfter("so beautiful.")

# --- Text Placeholder 18 (shape id 15): Client / Team block -----------------
$clientBox = Get-ShapeById $s 15
$clientRange = $clientBox.TextFrame.TextRange

$fullText = $clientRange.Text
$oldClientLabel = "Client(s):"
$newClientLabel = "Client:"
$idx = $fullText.IndexOf($oldClientLabel)
if ($idx -ge 0) {
    $clientRange.Characters($idx + 1, $oldClientLabel.Length).Text = $newClientLabel
}

$fullText = $clientRange.Text
$oldClientName = "Kim Whitehall, Nasa JPL"
$newClientName = "Kim Whitehall"
$idx = $fullText.IndexOf($oldClientName)
if ($idx -ge 0) {
    $clientRange.Characters($idx + 1, $oldClientName.Length).Text = $newClientName
}

# --- Title 1 (shape id 16): poster / team number -----------------------------
$numberBox = Get-ShapeById $s 16
$numberBox.TextFrame.TextRange.Text = "CS19"
